# LOM3070.xlsx update
# - Insert two rows for "Docentes responsáveis:" (professor names) right
#   under the "Objetivos:" block.
# - Fill in real Portuguese content for Objetivos / Programa resumido /
#   Programa / Bibliografia answer cells that previously held placeholder
#   (wrong) values copied from other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank rows at 13-14 (shifts everything from the old
#    row 13 "Programa resumido:" downward by two rows).
$ws.Rows("13:14").Insert()

# 2) The inserted rows inherit row 12's (single-column) formatting; copy
#    the B/C answer-column formatting from row 15 ("Programa resumido:",
#    which used to be row 13 before the insert) onto the two new rows,
#    then drop the stray column-A cells so rows 13/14 only have B & C,
#    matching the "Docentes responsáveis:" sub-rows.
$ws.Range("B15:C15").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$ws.Range("A13:A14").Clear()

# 3) Names of the two responsible faculty members.
$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("B14").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C14").Value = "519033 - Carlos Yujiro Shigue"

# 4) Objetivos: (row 10) had the wrong text in B/C - fix it.
$objetivos = "Fornecer oportunidade de aplicação dos conhecimentos de Engenharia de Materiais em empresa ou instituição de pesquisa sob a supervisão de docente do Departamento de Engenharia de Materiais da EEL. Complementação da formação geral curricular. Integração de conhecimentos técnicos, econômicos e de gestão. Adaptação psicológica e social do estudante à sua futura atividade profissional."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# 5) Programa resumido: (now row 15, after the insert) - fix placeholder.
$programaResumido = "Obtenção de estágio. Elaboração do plano de trabalho de estágio. Realização do estágio. Elaboração de relatórios parciais e/ou final."
$ws.Range("B15").Value = $programaResumido
$ws.Range("C15").Value = $programaResumido

# 6) Programa: (now row 17) - fix placeholder.
$programa = "Participação do aluno em processo seletivo de empresas, instituições de pesquisa ou no setor acadêmico. O estágio realizado sob a supervisão de docente designado pelo Coordenador de Estágio do curso de Engenharia de Materiais. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor, desde que relacionado com as áreas afins da Engenharia de Materiais. Apresentação de relatórios parciais e/ou final sobre as atividades desenvolvidas no estágio."
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa

# 7) Método: (now row 20) - fix placeholder (was showing the wrong text).
$metodo = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# 8) Critério: (now row 21) - fix placeholder.
$criterio = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# 9) Norma de recuperação: (now row 22) - fix placeholder.
$normaRecuperacao = "Não será oferecida recuperação."
$ws.Range("B22").Value = $normaRecuperacao
$ws.Range("C22").Value = $normaRecuperacao

# 10) Bibliografia: (now row 23) - previously empty of its own text, add it.
$bibliografia = "A ser definida com o supervisor responsável pelo estágio e pelo docente orientador em função das atividades desenvolvidas no estágio."
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia
